# chore: update Sheets via scheduled runner
# Applies numeric corrections to several leve-profit rows across sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per currentAveragePrice/
# LevePrice/LeveProfit recalculation.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2651.2
$ws.Cells.Item(137, 9).Value = 3410.4
$ws.Cells.Item(137, 11).Value = 10231.2
$ws.Cells.Item(137, 13).Value = -7681.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 36877
$ws.Cells.Item(7, 9).Value = 43673.5
$ws.Cells.Item(7, 10).Value = 32346
$ws.Cells.Item(7, 11).Value = 43673.5
$ws.Cells.Item(7, 12).Value = 32346
$ws.Cells.Item(7, 13).Value = -43559.5
$ws.Cells.Item(7, 14).Value = -32574
$ws.Cells.Item(61, 8).Value = 4291.5
$ws.Cells.Item(61, 9).Value = 3324.258
$ws.Cells.Item(61, 11).Value = 3324.258
$ws.Cells.Item(61, 13).Value = -3112.258
$ws.Cells.Item(74, 8).Value = 1642.8572
$ws.Cells.Item(74, 9).Value = 1101
$ws.Cells.Item(74, 11).Value = 1101
$ws.Cells.Item(74, 13).Value = -227
$ws.Cells.Item(77, 8).Value = 1642.8572
$ws.Cells.Item(77, 9).Value = 1101
$ws.Cells.Item(77, 11).Value = 5505
$ws.Cells.Item(77, 13).Value = -1137
$ws.Cells.Item(120, 8).Value = 40380.168
$ws.Cells.Item(120, 9).Value = 40381
$ws.Cells.Item(120, 10).Value = 40380
$ws.Cells.Item(120, 11).Value = 40381
$ws.Cells.Item(120, 12).Value = 40380
$ws.Cells.Item(120, 13).Value = -35543
$ws.Cells.Item(120, 14).Value = -50056
$ws.Cells.Item(122, 8).Value = 1652.1
$ws.Cells.Item(122, 9).Value = 1652.1
$ws.Cells.Item(122, 11).Value = 4956.299999999999
$ws.Cells.Item(122, 13).Value = -2506.299999999999
$ws.Cells.Item(136, 8).Value = 4291.5
$ws.Cells.Item(136, 9).Value = 3324.258
$ws.Cells.Item(136, 11).Value = 9972.774
$ws.Cells.Item(136, 13).Value = -7422.773999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value = 2887
$ws.Cells.Item(106, 10).Value = 2887
$ws.Cells.Item(106, 12).Value = 2887
$ws.Cells.Item(106, 14).Value = -5411
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1844.4
$ws.Cells.Item(31, 9).Value = 1564.5714
$ws.Cells.Item(31, 11).Value = 1564.5714
$ws.Cells.Item(31, 13).Value = -1269.5714
$ws.Cells.Item(34, 8).Value = 1844.4
$ws.Cells.Item(34, 9).Value = 1564.5714
$ws.Cells.Item(34, 11).Value = 1564.5714
$ws.Cells.Item(34, 13).Value = -1362.5714
$ws.Cells.Item(63, 8).Value = 81448.86
$ws.Cells.Item(63, 10).Value = 81448.86
$ws.Cells.Item(63, 12).Value = 81448.86
$ws.Cells.Item(63, 14).Value = -82820.86
$ws.Cells.Item(66, 8).Value = 81448.86
$ws.Cells.Item(66, 10).Value = 81448.86
$ws.Cells.Item(66, 12).Value = 244346.58
$ws.Cells.Item(66, 14).Value = -251210.58
$ws.Cells.Item(134, 8).Value = 3089.52
$ws.Cells.Item(134, 9).Value = 3105.1667
$ws.Cells.Item(134, 11).Value = 9315.500100000001
$ws.Cells.Item(134, 13).Value = -6780.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 48999.137
$ws.Cells.Item(34, 9).Value = 996.3333
$ws.Cells.Item(34, 10).Value = 56578.527
$ws.Cells.Item(34, 11).Value = 2988.9999
$ws.Cells.Item(34, 12).Value = 169735.581
$ws.Cells.Item(34, 13).Value = -2904.9999
$ws.Cells.Item(34, 14).Value = -169903.581
$ws.Cells.Item(131, 8).Value = 781.7
$ws.Cells.Item(131, 9).Value = 646.3333
$ws.Cells.Item(131, 11).Value = 1938.9999
$ws.Cells.Item(131, 13).Value = 3101.0001
$ws.Cells.Item(132, 8).Value = 1461.5
$ws.Cells.Item(132, 10).Value = 1360
$ws.Cells.Item(132, 12).Value = 12240
$ws.Cells.Item(132, 14).Value = -17300

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 1463
$ws.Cells.Item(31, 9).Value = 944.5
$ws.Cells.Item(31, 11).Value = 944.5
$ws.Cells.Item(31, 13).Value = -652.5
$ws.Cells.Item(37, 8).Value = 1463
$ws.Cells.Item(37, 9).Value = 944.5
$ws.Cells.Item(37, 11).Value = 944.5
$ws.Cells.Item(37, 13).Value = -667.5
$ws.Cells.Item(108, 8).Value = 98500
$ws.Cells.Item(108, 10).Value = 98500
$ws.Cells.Item(108, 12).Value = 98500
$ws.Cells.Item(108, 14).Value = -106180
$ws.Cells.Item(113, 8).Value = 964.8333
$ws.Cells.Item(113, 9).Value = 847.55554
$ws.Cells.Item(113, 10).Value = 1316.6666
$ws.Cells.Item(113, 11).Value = 847.55554
$ws.Cells.Item(113, 12).Value = 1316.6666
$ws.Cells.Item(113, 13).Value = 1322.44446
$ws.Cells.Item(113, 14).Value = -5656.6666
$ws.Cells.Item(122, 8).Value = 2495.4
$ws.Cells.Item(122, 9).Value = 1360.4445
$ws.Cells.Item(122, 10).Value = 4197.8335
$ws.Cells.Item(122, 11).Value = 4081.3335
$ws.Cells.Item(122, 12).Value = 12593.5005
$ws.Cells.Item(122, 13).Value = -1631.3335
$ws.Cells.Item(122, 14).Value = -17493.5005
$ws.Cells.Item(132, 8).Value = 2791.125
$ws.Cells.Item(132, 9).Value = 2917
$ws.Cells.Item(132, 10).Value = 2413.5
$ws.Cells.Item(132, 11).Value = 8751
$ws.Cells.Item(132, 12).Value = 7240.5
$ws.Cells.Item(132, 13).Value = -6221
$ws.Cells.Item(132, 14).Value = -12300.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 749
$ws.Cells.Item(22, 9).Value = 749
$ws.Cells.Item(22, 11).Value = 749
$ws.Cells.Item(22, 13).Value = -454
$ws.Cells.Item(27, 8).Value = 749
$ws.Cells.Item(27, 9).Value = 749
$ws.Cells.Item(27, 11).Value = 749
$ws.Cells.Item(27, 13).Value = -642
$ws.Cells.Item(68, 8).Value = 4446.2
$ws.Cells.Item(68, 9).Value = 4295.6665
$ws.Cells.Item(68, 11).Value = 4295.6665
$ws.Cells.Item(68, 13).Value = -3546.6665
$ws.Cells.Item(71, 8).Value = 4446.2
$ws.Cells.Item(71, 9).Value = 4295.6665
$ws.Cells.Item(71, 11).Value = 21478.3325
$ws.Cells.Item(71, 13).Value = -17734.3325
$ws.Cells.Item(95, 8).Value = 44744.75
$ws.Cells.Item(95, 10).Value = 44744.75
$ws.Cells.Item(95, 12).Value = 44744.75
$ws.Cells.Item(95, 14).Value = -50236.75
$ws.Cells.Item(99, 8).Value = 24999.2
$ws.Cells.Item(132, 8).Value = 2593.889
$ws.Cells.Item(132, 9).Value = 2625
$ws.Cells.Item(132, 10).Value = 2345
$ws.Cells.Item(132, 11).Value = 7875
$ws.Cells.Item(132, 12).Value = 7035
$ws.Cells.Item(132, 13).Value = -5345
$ws.Cells.Item(132, 14).Value = -12095
$ws.Cells.Item(136, 8).Value = 3364.8
$ws.Cells.Item(136, 9).Value = 2248
$ws.Cells.Item(136, 11).Value = 6744
$ws.Cells.Item(136, 13).Value = -4194

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4090
$ws.Cells.Item(81, 9).Value = 4143.5
$ws.Cells.Item(81, 11).Value = 8287
$ws.Cells.Item(81, 13).Value = -7226
$ws.Cells.Item(84, 8).Value = 4090
$ws.Cells.Item(84, 9).Value = 4143.5
$ws.Cells.Item(84, 11).Value = 41435
$ws.Cells.Item(84, 13).Value = -36131
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 11056.839
$ws.Cells.Item(136, 9).Value = 11630.069
$ws.Cells.Item(136, 11).Value = 34890.20699999999
$ws.Cells.Item(136, 13).Value = -32340.20699999999
